$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.992.92"
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = "'1.857.13"
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'311.59"
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').Value = "'1.003"
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').Value = "'0.5086"
$ws.Range('E7').Value = '  +1.92%  '
$ws.Range('E8').Value = '  -0.65%  '
$ws.Range('D9').Value = "'0.08256"
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').Value = "'41.48"
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = "'6.193"
$ws.Range('E12').Value = '  -2.64%  '
$ws.Range('D13').Value = "'20.53"
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = "'1.859.84"
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = "'7.191"
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').Value = "'0.00001097"
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = "'90.48"
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').Value = "'0.06598"
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('D20').Value = "'17.66"
$ws.Range('E20').Value = '  -1.95%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = "'6.013"
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('D23').Value = "'28.001.86"
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = "'11.04"
$ws.Range('E24').Value = '  -3.84%  '
$ws.Range('D25').Value = "'2.243"
$ws.Range('E25').Value = '  -1.66%  '
$ws.Range('D26').Value = "'2.546"
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('D27').Value = "'2.069.16"
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').Value = "'157.89"
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').Value = "'20.41"
$ws.Range('D30').Value = "'124.36"
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').Value = "'1.038"
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('D33').Value = "'5.613"
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').Value = "'3.605"
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').Value = "'9.456"
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').Value = "'0.06529"
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('D37').Value = "'0.02409"
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').Value = "'0.2166"
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').Value = "'1.202"
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('D40').Value = "'0.6405"
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').Value = "'1.236"
$ws.Range('E41').Value = '  -3.76%  '
$ws.Range('D42').Value = "'4.864"
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('D43').Value = "'11.15"
$ws.Range('E43').Value = '  -3.49%  '
$ws.Range('D44').Value = "'0.6055"
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('D45').Value = "'13.08"
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').Value = "'1.277"
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').Value = "'3.647"
$ws.Range('E47').Value = '  -0.78%  '
$ws.Range('D48').Value = "'1.993"
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').Value = "'1.207"
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('D50').Value = "'119.89"
$ws.Range('E50').Value = '  -0.84%  '
$ws.Range('D51').Value = "'78.79"
$ws.Range('E51').Value = '  -0.02%  '
